$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$old = "dnasr281@gmail.com, System"
$new = "System, dnasr281@gmail.com"

# Use Find/FindNext so that only cells that actually contain the exact
# old value are touched (avoids disturbing any other, untouched cells).
$first = $ws.Cells.Find($old)
if ($first -ne $null) {
    $firstAddress = $first.Address()
    $current = $first
    $continue = $true
    while ($continue) {
        $current.Value2 = $new
        $current = $ws.Cells.FindNext($current)
        if ($current -eq $null -or $current.Address() -eq $firstAddress) {
            $continue = $false
        }
    }
}
